# The "大房间" (Big Room) card's effect text has been rewritten as part of a
# rules overhaul: instead of doubling the hand-size bonus, the card now
# grants extra draws on reshuffle and extra reveals/picks when flipping.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "持续：房间宽度加1，高度加1，玩家重整时可以额外抽2张牌，翻选时可以额外翻2张牌、额外选1张牌。"

# The active selection moved from C8 to C4 in the saved view state.
$ws.Range("C4").Select()
